$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data table (rows 1-13, cols A-E) before rebuilding it,
# so leftover cells (e.g. old rows 12-13) don't linger.
$ws.Range("A1:F13").Clear()

# --- Column A: equipment list (unchanged) ---
$ws.Range("A1").Value = "Equipment Used"
$ws.Range("A2").Value = "Force Sensor - SHIMPO 30kg"
$ws.Range("A3").Value = "Multimeter - FLUKE 117"

# --- Headers ---
$ws.Range("C1").Value = "Sensor 1"
$ws.Range("E1").Value = "Sensor 2"

$ws.Range("C2").Value = "Applied Force [kg]"
$ws.Range("D2").Value = "Resistance [Ohms]"
$ws.Range("E2").Value = "Applied Force [kg]"
$ws.Range("F2").Value = "Resistance [Ohms]"

# --- Sensor 1 data (C/D) ---
$ws.Range("C3").Value = 0.2
$ws.Range("C4").Value = 0.4
$ws.Range("C5").Value = 0.6
$ws.Range("C6").Value = 0.8
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 2
$ws.Range("C9").Value = 3
$ws.Range("C10").Value = 4
$ws.Range("C11").Value = 5

$ws.Range("D3").Value = 150000
$ws.Range("D4").Value = 30000
$ws.Range("D5").Value = 8000
$ws.Range("D6").Value = 7900
$ws.Range("D7").Value = 6000
$ws.Range("D8").Value = 2000
$ws.Range("D9").Formula = "=0.5*1000"
$ws.Range("D10").Value = 450
$ws.Range("D11").Value = 230

# --- Sensor 2 data (E/F) ---
$ws.Range("E3").Value = 0.2
$ws.Range("E4").Value = 0.4
$ws.Range("E5").Value = 0.6
$ws.Range("E6").Value = 0.8
$ws.Range("E7").Value = 1
$ws.Range("E8").Value = 2
$ws.Range("E9").Value = 3
$ws.Range("E10").Value = 4
$ws.Range("E11").Value = 5

$ws.Range("F3").Value = "N/A"
$ws.Range("F4").Value = 24000
$ws.Range("F5").Value = 10000
$ws.Range("F6").Value = 6000
$ws.Range("F7").Value = 5900
$ws.Range("F8").Value = 1500
$ws.Range("F9").Value = 1000
$ws.Range("F10").Value = 600
$ws.Range("F11").Value = "N/A"

# --- Column widths to match the new (wider) header text ---
$ws.Columns.Item(3).ColumnWidth = 14.6
$ws.Columns.Item(4).ColumnWidth = 14.92
$ws.Columns.Item(5).ColumnWidth = 14.6
$ws.Columns.Item(6).ColumnWidth = 14.92

# --- Selection matches the authored state ---
$ws.Range("F3").Select()
